$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 294
$ws.Range("I33").Value = 220.14285
$ws.Range("J33").Value = 466.33334
$ws.Range("K33").Value = 220.14285
$ws.Range("L33").Value = 466.33334
$ws.Range("M33").Value = 8.85714999999999
$ws.Range("N33").Value = -924.33334

$ws.Range("H100").Value = 16668387
$ws.Range("I100").Value = 18183514
$ws.Range("K100").Value = 18183514
$ws.Range("M100").Value = -18182973

$ws.Range("H112").Value = 535587.5600000001
$ws.Range("I112").Value = 685
$ws.Range("J112").Value = 567052.4399999999
$ws.Range("K112").Value = 2055
$ws.Range("L112").Value = 1701157.32
$ws.Range("M112").Value = -947
$ws.Range("N112").Value = -1703373.32

$ws.Range("H129").Value = 822.8
$ws.Range("J129").Value = 899.81396
$ws.Range("L129").Value = 2699.44188
$ws.Range("N129").Value = -12699.44188

$ws.Range("H137").Value = 2723.08
$ws.Range("I137").Value = 1768.625
$ws.Range("J137").Value = 4419.8887
$ws.Range("K137").Value = 5305.875
$ws.Range("L137").Value = 13259.6661
$ws.Range("M137").Value = -2755.875
$ws.Range("N137").Value = -18359.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5086.4116
$ws.Range("I32").Value = 4518.548
$ws.Range("J32").Value = 7736.4443
$ws.Range("K32").Value = 4518.548
$ws.Range("L32").Value = 7736.4443
$ws.Range("M32").Value = -4231.548
$ws.Range("N32").Value = -8310.444299999999

$ws.Range("H74").Value = 5923.55
$ws.Range("I74").Value = 8544.5
$ws.Range("J74").Value = 3302.6
$ws.Range("K74").Value = 8544.5
$ws.Range("L74").Value = 3302.6
$ws.Range("M74").Value = -7670.5
$ws.Range("N74").Value = -5050.6

$ws.Range("H77").Value = 5923.55
$ws.Range("I77").Value = 8544.5
$ws.Range("J77").Value = 3302.6
$ws.Range("K77").Value = 42722.5
$ws.Range("L77").Value = 16513
$ws.Range("M77").Value = -38354.5
$ws.Range("N77").Value = -25249

$ws.Range("H101").Value = 42801
$ws.Range("J101").Value = 42801
$ws.Range("L101").Value = 42801
$ws.Range("N101").Value = -49291

$ws.Range("H102").Value = 2069
$ws.Range("I102").Value = 1953.5
$ws.Range("K102").Value = 1953.5
$ws.Range("M102").Value = -331.5

$ws.Range("H132").Value = 1891.7646
$ws.Range("I132").Value = 1160.5946
$ws.Range("K132").Value = 3481.7838
$ws.Range("M132").Value = -951.7837999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1291.4445
$ws.Range("I94").Value = 1216.6
$ws.Range("J94").Value = 1385
$ws.Range("K94").Value = 1216.6
$ws.Range("L94").Value = 1385
$ws.Range("M94").Value = -765.5999999999999
$ws.Range("N94").Value = -2287

$ws.Range("H103").Value = 35079.58
$ws.Range("J103").Value = 35079.58
$ws.Range("L103").Value = 35079.58
$ws.Range("N103").Value = -37423.58

$ws.Range("H107").Value = 2037.95
$ws.Range("I107").Value = 1897.5883
$ws.Range("K107").Value = 1897.5883
$ws.Range("M107").Value = 22.41170000000011

$ws.Range("H132").Value = 50750
$ws.Range("J132").Value = 50750
$ws.Range("L132").Value = 50750
$ws.Range("N132").Value = -60870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4347.933
$ws.Range("I31").Value = 1411.375
$ws.Range("K31").Value = 1411.375
$ws.Range("M31").Value = -1116.375

$ws.Range("H34").Value = 4347.933
$ws.Range("I34").Value = 1411.375
$ws.Range("K34").Value = 1411.375
$ws.Range("M34").Value = -1209.375

$ws.Range("H38").Value = 20999.938
$ws.Range("J38").Value = 20999.938
$ws.Range("L38").Value = 20999.938
$ws.Range("N38").Value = -21753.938

$ws.Range("H46").Value = 20999.938
$ws.Range("J46").Value = 20999.938
$ws.Range("L46").Value = 20999.938
$ws.Range("N46").Value = -21421.938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1398.5
$ws.Range("I5").Value = 379.75
$ws.Range("J5").Value = 2213.5
$ws.Range("K5").Value = 1139.25
$ws.Range("L5").Value = 6640.5
$ws.Range("M5").Value = -1027.25
$ws.Range("N5").Value = -6864.5

$ws.Range("H107").Value = 59257.766
$ws.Range("I107").Value = 431.75
$ws.Range("J107").Value = 200440.2
$ws.Range("K107").Value = 1295.25
$ws.Range("L107").Value = 601320.6000000001
$ws.Range("M107").Value = 624.75
$ws.Range("N107").Value = -605160.6000000001

$ws.Range("H113").Value = 599.2361
$ws.Range("I113").Value = 593.3333
$ws.Range("K113").Value = 1779.9999
$ws.Range("M113").Value = 390.0001

$ws.Range("H120").Value = 1950.25
$ws.Range("J120").Value = 5000
$ws.Range("L120").Value = 15000
$ws.Range("N120").Value = -24676

$ws.Range("H122").Value = 2622.932
$ws.Range("I122").Value = 872.4
$ws.Range("K122").Value = 7851.599999999999
$ws.Range("M122").Value = -5401.599999999999

$ws.Range("H132").Value = 3024
$ws.Range("I132").Value = 658.3333
$ws.Range("K132").Value = 5924.9997
$ws.Range("M132").Value = -3394.9997

$ws.Range("H135").Value = 1398.5
$ws.Range("I135").Value = 379.75
$ws.Range("J135").Value = 2213.5
$ws.Range("K135").Value = 3417.75
$ws.Range("L135").Value = 19921.5
$ws.Range("M135").Value = -882.75
$ws.Range("N135").Value = -24991.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3711.4167
$ws.Range("I132").Value = 1923.3334
$ws.Range("K132").Value = 5770.0002
$ws.Range("M132").Value = -3240.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1615.3684
$ws.Range("I46").Value = 2685.5
$ws.Range("K46").Value = 2685.5
$ws.Range("M46").Value = -2497.5

$ws.Range("H68").Value = 716.42
$ws.Range("I68").Value = 714.5657
$ws.Range("J68").Value = 900
$ws.Range("K68").Value = 714.5657
$ws.Range("L68").Value = 900
$ws.Range("M68").Value = 34.43430000000001
$ws.Range("N68").Value = -2398

$ws.Range("H71").Value = 716.42
$ws.Range("I71").Value = 714.5657
$ws.Range("J71").Value = 900
$ws.Range("K71").Value = 3572.8285
$ws.Range("L71").Value = 4500
$ws.Range("M71").Value = 171.1714999999999
$ws.Range("N71").Value = -11988

$ws.Range("H82").Value = 2098.1904
$ws.Range("I82").Value = 720.375
$ws.Range("K82").Value = 720.375
$ws.Range("M82").Value = -359.375

$ws.Range("H85").Value = 2098.1904
$ws.Range("I85").Value = 720.375
$ws.Range("K85").Value = 720.375
$ws.Range("M85").Value = 527.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 114800.1
$ws.Range("I62").Value = 5499.6665
$ws.Range("J62").Value = 278750.75
$ws.Range("K62").Value = 5499.6665
$ws.Range("L62").Value = 278750.75
$ws.Range("M62").Value = -4875.6665
$ws.Range("N62").Value = -279998.75

$ws.Range("H65").Value = 114800.1
$ws.Range("I65").Value = 5499.6665
$ws.Range("J65").Value = 278750.75
$ws.Range("K65").Value = 27498.3325
$ws.Range("L65").Value = 1393753.75
$ws.Range("M65").Value = -24378.3325
$ws.Range("N65").Value = -1399993.75
